$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("registerUsr")

# Update the two email values shown in column G (keep existing hyperlink formatting)
$ws.Range("G2").Value = "vijefg@gmail.com"
$ws.Range("G3").Value = "divyfg@gmail.com"

# Update the active selection on the sheet from L9 to G9
$ws.Activate()
$ws.Range("G9").Select()
